$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp from 10:16 to 10:46
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 10:46"

# Cantabria ("A35") deaths count (column E, "Muertes") increases from 3 to 4
$ws.Range("E35").Value = 4
